$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update comment text in column D that reference shared strings whose
# text content changed / got remapped.
$ws.Range("D10").Value = "AUROC + Matrice de confusion "
$ws.Range("D15").Value = "Matrice de confusion et AUROC"
$ws.Range("D23").Value = "Fait"

# Update progress values in column C
$ws.Range("C8").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 0.5

# Row 10 shrinks since the new comment text is shorter (matches the
# post-edit auto-computed wrapped height).
$ws.Rows("10:10").RowHeight = 49.5

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22:C22").Select()
